$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet, and add the new "franks" sheet right after it ---
$wsLawrence = $wb.Worksheets.Item(1)
$wsLawrence.Name = "UK_EQ5D_lawrence"

$wsFranks = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsLawrence)
$wsFranks.Name = "UK_EQ5D_franks"

# --- Match column widths on the new sheet to the lawrence sheet (A=15, B=17.7109375 chars) ---
$wsFranks.Columns.Item(1).ColumnWidth = 14.166666666666666
$wsFranks.Columns.Item(2).ColumnWidth = 16.877604166666668

# --- Header row (identical labels to the lawrence sheet) ---
$wsFranks.Cells.Item(1,1).Value = "REGRESSOR"
$wsFranks.Cells.Item(1,2).Value = "COEFFICIENT"
$wsFranks.Cells.Item(1,3).Value = "Constant"
$wsFranks.Cells.Item(1,4).Value = "Dhe_pcs"
$wsFranks.Cells.Item(1,5).Value = "Dhe_mcs"
$wsFranks.Cells.Item(1,6).Value = "Dhe_pcs_sq"
$wsFranks.Cells.Item(1,7).Value = "Dhe_mcs_sq"
$wsFranks.Cells.Item(1,8).Value = "Dhe_mcs_times_pcs"
$wsFranks.Cells.Item(1,9).Value = "Che_pcs_cb"

# --- Row 2: Constant ---
$wsFranks.Cells.Item(2,1).Value = "Constant"
$wsFranks.Cells.Item(2,2).Value = 0.8469
$wsFranks.Cells.Item(2,3).Value = 0
$wsFranks.Cells.Item(2,4).Value = 0
$wsFranks.Cells.Item(2,5).Value = 0
$wsFranks.Cells.Item(2,6).Value = 0
$wsFranks.Cells.Item(2,7).Value = 0
$wsFranks.Cells.Item(2,8).Value = 0
$wsFranks.Cells.Item(2,9).Value = 0

# --- Row 3: Dhe_pcs_c ---
$wsFranks.Cells.Item(3,1).Value = "Dhe_pcs_c"
$wsFranks.Cells.Item(3,2).Value = 0.01261
$wsFranks.Cells.Item(3,3).Value = 0
$wsFranks.Cells.Item(3,4).Value = 0
$wsFranks.Range("D3").NumberFormat = "0.00E+00"
$wsFranks.Cells.Item(3,5).Value = 0
$wsFranks.Cells.Item(3,6).Value = 0
$wsFranks.Cells.Item(3,7).Value = 0
$wsFranks.Cells.Item(3,8).Value = 0
$wsFranks.Cells.Item(3,9).Value = 0

# --- Row 4: Dhe_mcs_c ---
$wsFranks.Cells.Item(4,1).Value = "Dhe_mcs_c"
$wsFranks.Cells.Item(4,2).Value = 0.00759
$wsFranks.Cells.Item(4,3).Value = 0
$wsFranks.Cells.Item(4,4).Value = 0
$wsFranks.Cells.Item(4,5).Value = 0
$wsFranks.Range("E4").NumberFormat = "0.00E+00"
$wsFranks.Cells.Item(4,6).Value = 0
$wsFranks.Cells.Item(4,7).Value = 0
$wsFranks.Cells.Item(4,8).Value = 0
$wsFranks.Cells.Item(4,9).Value = 0

# --- Row 5: Dhe_pcs_c_sq ---
$wsFranks.Cells.Item(5,1).Value = "Dhe_pcs_c_sq"
$wsFranks.Cells.Item(5,2).Value = -0.00009
$wsFranks.Cells.Item(5,3).Value = 0
$wsFranks.Cells.Item(5,4).Value = 0
$wsFranks.Cells.Item(5,5).Value = 0
$wsFranks.Cells.Item(5,6).Value = 0
$wsFranks.Range("F5").NumberFormat = "0.00E+00"
$wsFranks.Cells.Item(5,7).Value = 0
$wsFranks.Cells.Item(5,8).Value = 0
$wsFranks.Cells.Item(5,9).Value = 0

# --- Row 6: Dhe_mcs_c_sq ---
$wsFranks.Cells.Item(6,1).Value = "Dhe_mcs_c_sq"
$wsFranks.Cells.Item(6,2).Value = -0.00015
$wsFranks.Cells.Item(6,3).Value = 0
$wsFranks.Cells.Item(6,4).Value = 0
$wsFranks.Cells.Item(6,5).Value = 0
$wsFranks.Cells.Item(6,6).Value = 0
$wsFranks.Cells.Item(6,7).Value = 0
$wsFranks.Range("G6").NumberFormat = "0.00E+00"
$wsFranks.Cells.Item(6,8).Value = 0
$wsFranks.Cells.Item(6,9).Value = 0

# --- Row 7: Dhe_mcs_c_times_pcs_c ---
$wsFranks.Cells.Item(7,1).Value = "Dhe_mcs_c_times_pcs_c"
$wsFranks.Cells.Item(7,2).Value = -0.00015
$wsFranks.Cells.Item(7,3).Value = 0
$wsFranks.Cells.Item(7,4).Value = 0
$wsFranks.Cells.Item(7,5).Value = 0
$wsFranks.Cells.Item(7,6).Value = 0
$wsFranks.Cells.Item(7,7).Value = 0
$wsFranks.Cells.Item(7,8).Value = 0
$wsFranks.Range("H7").NumberFormat = "0.00E+00"
$wsFranks.Cells.Item(7,9).Value = 0

# --- Row 8: Dhe_pcs_cb ---
$wsFranks.Cells.Item(8,1).Value = "Dhe_pcs_cb"
$wsFranks.Cells.Item(8,2).Value = 0
$wsFranks.Range("B8").NumberFormat = "0.00E+00"
$wsFranks.Cells.Item(8,3).Value = 0
$wsFranks.Cells.Item(8,4).Value = 0
$wsFranks.Cells.Item(8,5).Value = 0
$wsFranks.Cells.Item(8,6).Value = 0
$wsFranks.Cells.Item(8,7).Value = 0
$wsFranks.Cells.Item(8,8).Value = 0
$wsFranks.Cells.Item(8,9).Value = 0
$wsFranks.Range("I8").NumberFormat = "0.00E+00"

# --- Zero out the previously non-zero higher-order-term coefficients on the lawrence sheet ---
$wsLawrence.Range("C2").Value = 0
$wsLawrence.Range("D3").Value = 0
$wsLawrence.Range("E4").Value = 0
$wsLawrence.Range("F5").Value = 0
$wsLawrence.Range("G6").Value = 0
$wsLawrence.Range("H7").Value = 0
$wsLawrence.Range("I8").Value = 0

# --- Selections: lawrence now shows E4 selected (no longer the active tab); franks is active, A7 selected ---
$wsLawrence.Range("E4").Select()
$wsFranks.Range("A7").Select()
